$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19 (shifts existing rows 19-32 down to 20-33,
# preserving their data and formatting).
$ws.Rows(19).Insert()

# Populate the newly inserted row 19 with the new weekly price record.
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C19").Value = "Arica y Parinacota"
$ws.Range("D19").Value = 45033
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = 100112013
$ws.Range("G19").Value = "Alcachofa"
$ws.Range("H19").Value = "Madrigal"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 120
$ws.Range("K19").Value = 23000
$ws.Range("L19").Value = 24000
$ws.Range("M19").Value = 23500
$ws.Range("N19").Value = "$/caja 40 unidades"
$ws.Range("O19").Value = "Provincia de Limarí"
$ws.Range("P19").Value = 588
$ws.Range("Q19").Value = 40
$ws.Range("R19").Value = "Hortaliza"
